# Update cryptos list (Price and Volume(1h) columns) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.187.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.434.61"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.66"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.68%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.79%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.31"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.32%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.52%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.812.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.427.40"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.832"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.124.82"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.65"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.92"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.65%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.18%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.88%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.84"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +16.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.62"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +8.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.16"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.46%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.53"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "130.47"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +18.93%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.28"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.73%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.958.26"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.05%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.668.79"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.36%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.42"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.70%  "
